# PISA 2015 scores by subject
#
# 1) Strip all whitespace (incl. the stray leading/doubled non-breaking
#    spaces) from every country name in column B (rows 2-71) on each of
#    the three sheets: "Singapore" instead of " Singapore",
#    "HongKong" instead of " Hong Kong", "Switzerland" instead of
#    "  Switzerland", etc. Rank (col A) and score (col C) values are left
#    untouched.
# 2) Move the active tab / selection from Math!D11 to Reading!E8, so
#    "Reading" is the sheet that's on top and selected when the workbook
#    is reopened.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Math", "Science", "Reading")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    for ($r = 2; $r -le 71; $r++) {
        $cell = $ws.Cells.Item($r, 2)
        $clean = $cell.Value2 -replace '\s', ''
        $cell.Value = $clean
    }
}

# Clear the old Math selection and land on Reading!E8 as the active view.
$wsMath = $wb.Worksheets.Item("Math")
$wsMath.Activate()
$wsMath.Range("A1").Select()

$wsReading = $wb.Worksheets.Item("Reading")
$wsReading.Activate()
$wsReading.Range("E8").Select()
